# "subimos el último SPA": refresh last week's "uds. Objetivo semana pasada" (col R)
# with the real figures (previously placeholder 0s). This ripples into the
# "Tendencia Consumo" (col T = Uds. Vtas. reales (S) - uds. Objetivo (R), floored at 0)
# and "Pedido Final" (col U) columns that were snapshotted alongside it, plus the
# "Total_Unidades" summary cell (C209). Three SKUs (rows 43, 49, 101) drop out of the
# current order (fully covered/no longer needed) and are hidden accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (R, T, U) new values. Only columns that actually changed are listed.
# row 5: R5=1, T5=2, U5=4
$ws.Range("R5").Value = 1
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = 4
# row 7: R7=1
$ws.Range("R7").Value = 1
# row 9: R9=1
$ws.Range("R9").Value = 1
# row 10: R10=3
$ws.Range("R10").Value = 3
# row 12: R12=3, T12=0
$ws.Range("R12").Value = 3
$ws.Range("T12").Value = 0
# row 15: R15=1
$ws.Range("R15").Value = 1
# row 17: R17=1
$ws.Range("R17").Value = 1
# row 21: R21=2
$ws.Range("R21").Value = 2
# row 22: R22=1
$ws.Range("R22").Value = 1
# row 25: R25=1
$ws.Range("R25").Value = 1
# row 26: R26=1
$ws.Range("R26").Value = 1
# row 34: R34=1
$ws.Range("R34").Value = 1
# row 40: R40=4, T40=1
$ws.Range("R40").Value = 4
$ws.Range("T40").Value = 1
# row 43: R43=1, T43=2, U43=0
$ws.Range("R43").Value = 1
$ws.Range("T43").Value = 2
$ws.Range("U43").Value = 0
# row 49: R49=3, T49=0, U49=0
$ws.Range("R49").Value = 3
$ws.Range("T49").Value = 0
$ws.Range("U49").Value = 0
# row 52: R52=1
$ws.Range("R52").Value = 1
# row 53: R53=2
$ws.Range("R53").Value = 2
# row 54: R54=5
$ws.Range("R54").Value = 5
# row 61: R61=4
$ws.Range("R61").Value = 4
# row 66: R66=2
$ws.Range("R66").Value = 2
# row 73: R73=1
$ws.Range("R73").Value = 1
# row 74: R74=1, T74=0, U74=2
$ws.Range("R74").Value = 1
$ws.Range("T74").Value = 0
$ws.Range("U74").Value = 2
# row 91: R91=3
$ws.Range("R91").Value = 3
# row 95: R95=1, T95=0, U95=2
$ws.Range("R95").Value = 1
$ws.Range("T95").Value = 0
$ws.Range("U95").Value = 2
# row 96: R96=1, T96=1
$ws.Range("R96").Value = 1
$ws.Range("T96").Value = 1
# row 101: R101=2, T101=1, U101=0
$ws.Range("R101").Value = 2
$ws.Range("T101").Value = 1
$ws.Range("U101").Value = 0
# row 102: R102=4, T102=1
$ws.Range("R102").Value = 4
$ws.Range("T102").Value = 1
# row 108: R108=1, T108=4
$ws.Range("R108").Value = 1
$ws.Range("T108").Value = 4
# row 109: R109=2
$ws.Range("R109").Value = 2
# row 110: R110=2
$ws.Range("R110").Value = 2
# row 111: R111=6
$ws.Range("R111").Value = 6
# row 112: R112=1
$ws.Range("R112").Value = 1
# row 114: R114=1, T114=3
$ws.Range("R114").Value = 1
$ws.Range("T114").Value = 3
# row 115: R115=1, T115=0
$ws.Range("R115").Value = 1
$ws.Range("T115").Value = 0
# row 118: R118=11, T118=0
$ws.Range("R118").Value = 11
$ws.Range("T118").Value = 0
# row 119: R119=4, T119=20
$ws.Range("R119").Value = 4
$ws.Range("T119").Value = 20
# row 120: R120=4
$ws.Range("R120").Value = 4
# row 121: R121=7
$ws.Range("R121").Value = 7
# row 122: R122=5, T122=0, U122=2
$ws.Range("R122").Value = 5
$ws.Range("T122").Value = 0
$ws.Range("U122").Value = 2
# row 123: R123=3
$ws.Range("R123").Value = 3
# row 124: R124=1
$ws.Range("R124").Value = 1
# row 126: R126=6, T126=10
$ws.Range("R126").Value = 6
$ws.Range("T126").Value = 10
# row 129: R129=17, T129=0
$ws.Range("R129").Value = 17
$ws.Range("T129").Value = 0
# row 130: R130=2
$ws.Range("R130").Value = 2
# row 131: R131=1
$ws.Range("R131").Value = 1
# row 132: R132=1, T132=0
$ws.Range("R132").Value = 1
$ws.Range("T132").Value = 0
# row 133: R133=13
$ws.Range("R133").Value = 13
# row 134: R134=3
$ws.Range("R134").Value = 3
# row 139: R139=1, T139=0
$ws.Range("R139").Value = 1
$ws.Range("T139").Value = 0
# row 140: R140=1, T140=0
$ws.Range("R140").Value = 1
$ws.Range("T140").Value = 0
# row 145: R145=2, T145=3
$ws.Range("R145").Value = 2
$ws.Range("T145").Value = 3
# row 152: R152=2, T152=3
$ws.Range("R152").Value = 2
$ws.Range("T152").Value = 3
# row 153: R153=1, T153=2, U153=5
$ws.Range("R153").Value = 1
$ws.Range("T153").Value = 2
$ws.Range("U153").Value = 5
# row 157: R157=2
$ws.Range("R157").Value = 2
# row 158: R158=1
$ws.Range("R158").Value = 1
# row 161: R161=1
$ws.Range("R161").Value = 1
# row 162: R162=1
$ws.Range("R162").Value = 1
# row 167: R167=2, T167=0
$ws.Range("R167").Value = 2
$ws.Range("T167").Value = 0
# row 172: R172=1
$ws.Range("R172").Value = 1
# row 175: R175=2
$ws.Range("R175").Value = 2
# row 176: R176=1, T176=5
$ws.Range("R176").Value = 1
$ws.Range("T176").Value = 5
# row 180: R180=3, T180=6
$ws.Range("R180").Value = 3
$ws.Range("T180").Value = 6
# row 186: R186=3, T186=0, U186=1
$ws.Range("R186").Value = 3
$ws.Range("T186").Value = 0
$ws.Range("U186").Value = 1
# row 190: R190=1, T190=2
$ws.Range("R190").Value = 1
$ws.Range("T190").Value = 2
# row 191: R191=3, T191=0
$ws.Range("R191").Value = 3
$ws.Range("T191").Value = 0
# row 192: R192=10, T192=0
$ws.Range("R192").Value = 10
$ws.Range("T192").Value = 0
# row 193: R193=2, T193=0
$ws.Range("R193").Value = 2
$ws.Range("T193").Value = 0
# row 194: R194=1, T194=0
$ws.Range("R194").Value = 1
$ws.Range("T194").Value = 0
# row 195: R195=15
$ws.Range("R195").Value = 15
# row 196: R196=13
$ws.Range("R196").Value = 13
# row 197: R197=10, T197=0
$ws.Range("R197").Value = 10
$ws.Range("T197").Value = 0
# row 198: R198=2, T198=0
$ws.Range("R198").Value = 2
$ws.Range("T198").Value = 0
# row 199: R199=2, T199=2
$ws.Range("R199").Value = 2
$ws.Range("T199").Value = 2
# row 206: R206=2
$ws.Range("R206").Value = 2

# Summary metric "Total_Unidades" (B209 label / C209 value)
$ws.Range("C209").Value = 163

# SKUs fully satisfied this week -> hide their rows (matches author toggling row visibility)
$ws.Rows.Item(43).Hidden = $true
$ws.Rows.Item(49).Hidden = $true
$ws.Rows.Item(101).Hidden = $true
